# Fixed naive component forecaster bug - Presentation state 11.02.
# New data row inserted at the top of the error table (Q0..Q9 block),
# shifting prior rows' metrics down by one and dropping the oldest row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.1724578193461484; C = 0.39058239716261;   D = 0.3033305724894426; E = 0.550754548314803;  F = 0.5414156770869448; G = 15 }
    3  = @{ B = 0.2400528213269932; C = 0.3753181292658889;  D = 0.2348887549825461; E = 0.4846532316848265; F = 0.4369199514236829; G = 14 }
    4  = @{ B = 0.273599112058131;  C = 0.3783303919539526;  D = 0.2266048251956944; E = 0.4760302776039508; F = 0.4054553575916278; G = 13 }
    5  = @{ B = 0.323052878118673;  C = 0.3688046801716363;  D = 0.2639536099994781; E = 0.5137641579552608; F = 0.4172513277134101; G = 12 }
    6  = @{ B = 0.3213959399964313; C = 0.3527736814977633;  D = 0.2405711618554991; E = 0.4904805417705163; F = 0.3885915500499728; G = 11 }
    7  = @{ B = 0.3040077233811113; C = 0.3040077233811113;  D = 0.2525516835738419; E = 0.5025452055027905; F = 0.4218095774931774; G = 10 }
    8  = @{ B = 0.3195075457514495; C = 0.3344878761007612;  D = 0.2334815320987104; E = 0.4831992674856931; F = 0.3844749898822031; G = 9 }
    9  = @{ B = 0.3427370110204659; C = 0.3605818570417189;  D = 0.2023703037307001; E = 0.4498558699524772; F = 0.3114971130477678; G = 8 }
    10 = @{ B = 0.2987280035122604; C = 0.2987280035122604;  D = 0.1585806734357395; E = 0.3982218896993729; F = 0.2844280147574629; G = 7 }
    11 = @{ B = 0.2911926455816474; C = 0.3419929353797495;  D = 0.3451110422066679; E = 0.587461524022355;  F = 0.5589109611011353; G = 6 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
